$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 169.5951436666667
$ws.Range("H2").Value = 508.785431
$ws.Range("I2").Value = 0.2074259764082431
$ws.Range("J2").Value = 0.2074259764082431
$ws.Range("M2").Value = 13.07123266666667
$ws.Range("N2").Value = 39.213698
$ws.Range("O2").Value = 0.1468300556961012
$ws.Range("P2").Value = 0.1468300556961012
$ws.Range("Q2").Value = 2216.81758200376
$ws.Range("R2").Value = 19951.35823803384
$ws.Range("S2").Value = 0.03045636766884051
$ws.Range("T2").Value = 0.03045636766884051
$ws.Range("G3").Value = 169.5951436666667
$ws.Range("H3").Value = 508.785431
$ws.Range("I3").Value = 0.2074259764082431
$ws.Range("J3").Value = 0.2074259764082431
$ws.Range("O3").Value = 0.04211112362724598
$ws.Range("P3").Value = 0.04211112362724597
$ws.Range("Q3").Value = 635.7872631202143
$ws.Range("R3").Value = 5722.085368081929
$ws.Range("S3").Value = 0.008734940936029733
$ws.Range("T3").Value = 0.008734940936029731
$ws.Range("G4").Value = 169.5951436666667
$ws.Range("H4").Value = 508.785431
$ws.Range("I4").Value = 0.2074259764082431
$ws.Range("J4").Value = 0.2074259764082431
$ws.Range("M4").Value = 40.10810466666667
$ws.Range("N4").Value = 120.324314
$ws.Range("O4").Value = 0.4505370986999281
$ws.Range("P4").Value = 0.450537098699928
$ws.Range("Q4").Value = 6802.139773141038
$ws.Range("R4").Value = 61219.25795826934
$ws.Range("S4").Value = 0.09345309760596957
$ws.Range("T4").Value = 0.09345309760596957
$ws.Range("G5").Value = 169.5951436666667
$ws.Range("H5").Value = 508.785431
$ws.Range("I5").Value = 0.2074259764082431
$ws.Range("J5").Value = 0.2074259764082431
$ws.Range("M5").Value = 2.890218
$ws.Range("N5").Value = 8.670653999999999
$ws.Range("O5").Value = 0.03246601760797012
$ws.Range("P5").Value = 0.03246601760797011
$ws.Range("Q5").Value = 490.166936937986
$ws.Range("R5").Value = 4411.502432441874
$ws.Range("S5").Value = 0.006734295402420414
$ws.Range("T5").Value = 0.006734295402420413
$ws.Range("G6").Value = 169.5951436666667
$ws.Range("H6").Value = 508.785431
$ws.Range("I6").Value = 0.2074259764082431
$ws.Range("J6").Value = 0.2074259764082431
$ws.Range("M6").Value = 29.20445966666667
$ws.Range("N6").Value = 87.61337900000001
$ws.Range("O6").Value = 0.3280557043687546
$ws.Range("P6").Value = 0.3280557043687546
$ws.Range("Q6").Value = 4952.934532875706
$ws.Range("R6").Value = 44576.41079588135
$ws.Range("S6").Value = 0.06804727479498288
$ws.Range("T6").Value = 0.06804727479498286
$ws.Range("I7").Value = 0.4056457917095931
$ws.Range("J7").Value = 0.405645791709593
$ws.Range("M7").Value = 13.07123266666667
$ws.Range("N7").Value = 39.213698
$ws.Range("O7").Value = 0.1468300556961012
$ws.Range("P7").Value = 0.1468300556961012
$ws.Range("Q7").Value = 4335.2464271776
$ws.Range("R7").Value = 39017.2178445984
$ws.Range("S7").Value = 0.05956099418960863
$ws.Range("T7").Value = 0.05956099418960861
$ws.Range("I8").Value = 0.4056457917095931
$ws.Range("J8").Value = 0.405645791709593
$ws.Range("O8").Value = 0.04211112362724598
$ws.Range("P8").Value = 0.04211112362724597
$ws.Range("S8").Value = 0.01708220008355475
$ws.Range("T8").Value = 0.01708220008355474
$ws.Range("I9").Value = 0.4056457917095931
$ws.Range("J9").Value = 0.405645791709593
$ws.Range("M9").Value = 40.10810466666667
$ws.Range("N9").Value = 120.324314
$ws.Range("O9").Value = 0.4505370986999281
$ws.Range("P9").Value = 0.450537098699928
$ws.Range("Q9").Value = 13302.38102948351
$ws.Range("R9").Value = 119721.4292653516
$ws.Range("S9").Value = 0.1827584780966754
$ws.Range("T9").Value = 0.1827584780966754
$ws.Range("I10").Value = 0.4056457917095931
$ws.Range("J10").Value = 0.405645791709593
$ws.Range("M10").Value = 2.890218
$ws.Range("N10").Value = 8.670653999999999
$ws.Range("O10").Value = 0.03246601760797012
$ws.Range("P10").Value = 0.03246601760797011
$ws.Range("Q10").Value = 958.5788561638118
$ws.Range("R10").Value = 8627.209705474306
$ws.Range("S10").Value = 0.01316970341624263
$ws.Range("T10").Value = 0.01316970341624262
$ws.Range("I11").Value = 0.4056457917095931
$ws.Range("J11").Value = 0.405645791709593
$ws.Range("M11").Value = 29.20445966666667
$ws.Range("N11").Value = 87.61337900000001
$ws.Range("O11").Value = 0.3280557043687546
$ws.Range("P11").Value = 0.3280557043687546
$ws.Range("Q11").Value = 9686.04359330525
$ws.Range("R11").Value = 87174.39233974725
$ws.Range("S11").Value = 0.1330744159235117
$ws.Range("T11").Value = 0.1330744159235117
$ws.Range("G12").Value = 98.17454766666667
$ws.Range("H12").Value = 294.523643
$ws.Range("I12").Value = 0.1200739064098473
$ws.Range("J12").Value = 0.1200739064098473
$ws.Range("M12").Value = 13.07123266666667
$ws.Range("N12").Value = 39.213698
$ws.Range("O12").Value = 0.1468300556961012
$ws.Range("P12").Value = 0.1468300556961012
$ws.Range("Q12").Value = 1283.262354495757
$ws.Range("R12").Value = 11549.36119046181
$ws.Range("S12").Value = 0.01763045836580632
$ws.Range("T12").Value = 0.01763045836580632
$ws.Range("G13").Value = 98.17454766666667
$ws.Range("H13").Value = 294.523643
$ws.Range("I13").Value = 0.1200739064098473
$ws.Range("J13").Value = 0.1200739064098473
$ws.Range("O13").Value = 0.04211112362724598
$ws.Range("P13").Value = 0.04211112362724597
$ws.Range("Q13").Value = 368.0419475438263
$ws.Range("R13").Value = 3312.377527894437
$ws.Range("S13").Value = 0.005056447117231443
$ws.Range("T13").Value = 0.005056447117231441
$ws.Range("G14").Value = 98.17454766666667
$ws.Range("H14").Value = 294.523643
$ws.Range("I14").Value = 0.1200739064098473
$ws.Range("J14").Value = 0.1200739064098473
$ws.Range("M14").Value = 40.10810466666667
$ws.Range("N14").Value = 120.324314
$ws.Range("O14").Value = 0.4505370986999281
$ws.Range("P14").Value = 0.450537098699928
$ws.Range("Q14").Value = 3937.595033417323
$ws.Range("R14").Value = 35438.3553007559
$ws.Range("S14").Value = 0.0540977494234593
$ws.Range("T14").Value = 0.05409774942345928
$ws.Range("G15").Value = 98.17454766666667
$ws.Range("H15").Value = 294.523643
$ws.Range("I15").Value = 0.1200739064098473
$ws.Range("J15").Value = 0.1200739064098473
$ws.Range("M15").Value = 2.890218
$ws.Range("N15").Value = 8.670653999999999
$ws.Range("O15").Value = 0.03246601760797012
$ws.Range("P15").Value = 0.03246601760797011
$ws.Range("Q15").Value = 283.7458448080579
$ws.Range("R15").Value = 2553.712603272522
$ws.Range("S15").Value = 0.003898321559759858
$ws.Range("T15").Value = 0.003898321559759857
$ws.Range("G16").Value = 98.17454766666667
$ws.Range("H16").Value = 294.523643
$ws.Range("I16").Value = 0.1200739064098473
$ws.Range("J16").Value = 0.1200739064098473
$ws.Range("M16").Value = 29.20445966666667
$ws.Range("N16").Value = 87.61337900000001
$ws.Range("O16").Value = 0.3280557043687546
$ws.Range("P16").Value = 0.3280557043687546
$ws.Range("Q16").Value = 2867.134617624411
$ws.Range("R16").Value = 25804.2115586197
$ws.Range("S16").Value = 0.03939092994359038
$ws.Range("T16").Value = 0.03939092994359036
$ws.Range("G17").Value = 105.935201
$ws.Range("H17").Value = 317.805603
$ws.Range("I17").Value = 0.1295656940897851
$ws.Range("J17").Value = 0.1295656940897851
$ws.Range("M17").Value = 13.07123266666667
$ws.Range("N17").Value = 39.213698
$ws.Range("O17").Value = 0.1468300556961012
$ws.Range("P17").Value = 0.1468300556961012
$ws.Range("Q17").Value = 1384.703659861099
$ws.Range("R17").Value = 12462.33293874989
$ws.Range("S17").Value = 0.01902413807950716
$ws.Range("T17").Value = 0.01902413807950716
$ws.Range("G18").Value = 105.935201
$ws.Range("H18").Value = 317.805603
$ws.Range("I18").Value = 0.1295656940897851
$ws.Range("J18").Value = 0.1295656940897851
$ws.Range("O18").Value = 0.04211112362724598
$ws.Range("P18").Value = 0.04211112362724597
$ws.Range("Q18").Value = 397.135496074453
$ws.Range("R18").Value = 3574.219464670077
$ws.Range("S18").Value = 0.005456156961664876
$ws.Range("T18").Value = 0.005456156961664873
$ws.Range("G19").Value = 105.935201
$ws.Range("H19").Value = 317.805603
$ws.Range("I19").Value = 0.1295656940897851
$ws.Range("J19").Value = 0.1295656940897851
$ws.Range("M19").Value = 40.10810466666667
$ws.Range("N19").Value = 120.324314
$ws.Range("O19").Value = 0.4505370986999281
$ws.Range("P19").Value = 0.450537098699928
$ws.Range("Q19").Value = 4248.860129592372
$ws.Range("R19").Value = 38239.74116633135
$ws.Range("S19").Value = 0.05837415190625422
$ws.Range("T19").Value = 0.0583741519062542
$ws.Range("G20").Value = 105.935201
$ws.Range("H20").Value = 317.805603
$ws.Range("I20").Value = 0.1295656940897851
$ws.Range("J20").Value = 0.1295656940897851
$ws.Range("M20").Value = 2.890218
$ws.Range("N20").Value = 8.670653999999999
$ws.Range("O20").Value = 0.03246601760797012
$ws.Range("P20").Value = 0.03246601760797011
$ws.Range("Q20").Value = 306.175824763818
$ws.Range("R20").Value = 2755.582422874362
$ws.Range("S20").Value = 0.004206482105707834
$ws.Range("T20").Value = 0.004206482105707832
$ws.Range("G21").Value = 105.935201
$ws.Range("H21").Value = 317.805603
$ws.Range("I21").Value = 0.1295656940897851
$ws.Range("J21").Value = 0.1295656940897851
$ws.Range("M21").Value = 29.20445966666667
$ws.Range("N21").Value = 87.61337900000001
$ws.Range("O21").Value = 0.3280557043687546
$ws.Range("P21").Value = 0.3280557043687546
$ws.Range("Q21").Value = 3093.780304884727
$ws.Range("R21").Value = 27844.02274396254
$ws.Range("S21").Value = 0.04250476503665105
$ws.Range("T21").Value = 0.04250476503665104
$ws.Range("G22").Value = 112.249611
$ws.Range("H22").Value = 336.748833
$ws.Range("I22").Value = 0.1372886313825315
$ws.Range("J22").Value = 0.1372886313825314
$ws.Range("M22").Value = 13.07123266666667
$ws.Range("N22").Value = 39.213698
$ws.Range("O22").Value = 0.1468300556961012
$ws.Range("P22").Value = 0.1468300556961012
$ws.Range("Q22").Value = 1467.240782123826
$ws.Range("R22").Value = 13205.16703911443
$ws.Range("S22").Value = 0.02015809739233861
$ws.Range("T22").Value = 0.0201580973923386
$ws.Range("G23").Value = 112.249611
$ws.Range("H23").Value = 336.748833
$ws.Range("I23").Value = 0.1372886313825315
$ws.Range("J23").Value = 0.1372886313825314
$ws.Range("O23").Value = 0.04211112362724598
$ws.Range("P23").Value = 0.04211112362724597
$ws.Range("Q23").Value = 420.807290946183
$ws.Range("R23").Value = 3787.265618515647
$ws.Range("S23").Value = 0.005781378528765185
$ws.Range("T23").Value = 0.005781378528765183
$ws.Range("G24").Value = 112.249611
$ws.Range("H24").Value = 336.748833
$ws.Range("I24").Value = 0.1372886313825315
$ws.Range("J24").Value = 0.1372886313825314
$ws.Range("M24").Value = 40.10810466666667
$ws.Range("N24").Value = 120.324314
$ws.Range("O24").Value = 0.4505370986999281
$ws.Range("P24").Value = 0.450537098699928
$ws.Range("Q24").Value = 4502.119146780618
$ws.Range("R24").Value = 40519.07232102557
$ws.Range("S24").Value = 0.06185362166756963
$ws.Range("T24").Value = 0.06185362166756961
$ws.Range("G25").Value = 112.249611
$ws.Range("H25").Value = 336.748833
$ws.Range("I25").Value = 0.1372886313825315
$ws.Range("J25").Value = 0.1372886313825314
$ws.Range("M25").Value = 2.890218
$ws.Range("N25").Value = 8.670653999999999
$ws.Range("O25").Value = 0.03246601760797012
$ws.Range("P25").Value = 0.03246601760797011
$ws.Range("Q25").Value = 324.425846205198
$ws.Range("R25").Value = 2919.832615846782
$ws.Range("S25").Value = 0.004457215123839385
$ws.Range("T25").Value = 0.004457215123839384
$ws.Range("G26").Value = 112.249611
$ws.Range("H26").Value = 336.748833
$ws.Range("I26").Value = 0.1372886313825315
$ws.Range("J26").Value = 0.1372886313825314
$ws.Range("M26").Value = 29.20445966666667
$ws.Range("N26").Value = 87.61337900000001
$ws.Range("O26").Value = 0.3280557043687546
$ws.Range("P26").Value = 0.3280557043687546
$ws.Range("Q26").Value = 3278.189237048523
$ws.Range("R26").Value = 29503.70313343671
$ws.Range("S26").Value = 0.04503831867001867
$ws.Range("T26").Value = 0.04503831867001866
